# Auto-generated edit script: refresh market-price derived columns (H-N)
# across multiple Leve-profit sheets, per scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 337.08334
$ws.Range("H40").Value = 1083.4
$ws.Range("J40").Value = 997.5
$ws.Range("L40").Value = 997.5
$ws.Range("N40").Value = -1347.5
$ws.Range("H53").Value = 578.9091
$ws.Range("I53").Value = 436.8
$ws.Range("K53").Value = 436.8
$ws.Range("M53").Value = 200.2
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H116").Value = 10510.625
$ws.Range("I116").Value = 3648.75
$ws.Range("J116").Value = 17372.5
$ws.Range("K116").Value = 3648.75
$ws.Range("L116").Value = 17372.5
$ws.Range("M116").Value = -206.75
$ws.Range("N116").Value = -24256.5
$ws.Range("H137").Value = 1421.2222
$ws.Range("I137").Value = 1311.375
$ws.Range("J137").Value = 2300
$ws.Range("K137").Value = 3934.125
$ws.Range("L137").Value = 6900
$ws.Range("M137").Value = -1384.125
$ws.Range("N137").Value = -12000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 48000
$ws.Range("J132").Value = 48000
$ws.Range("L132").Value = 48000
$ws.Range("N132").Value = -58120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2431.5
$ws.Range("J132").Value = 2369.3333
$ws.Range("L132").Value = 7107.999899999999
$ws.Range("N132").Value = -12167.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 218.61539
$ws.Range("J23").Value = 253
$ws.Range("L23").Value = 759
$ws.Range("N23").Value = -1229
$ws.Range("H34").Value = 804.1667
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 1005.55554
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 3016.66662
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -3184.66662
$ws.Range("H38").Value = 1421.0769
$ws.Range("I38").Value = 1136.1428
$ws.Range("J38").Value = 1753.5
$ws.Range("K38").Value = 3408.4284
$ws.Range("L38").Value = 5260.5
$ws.Range("M38").Value = -3061.4284
$ws.Range("N38").Value = -5954.5
$ws.Range("H39").Value = 1142.8572
$ws.Range("H55").Value = 775
$ws.Range("J55").Value = 1250
$ws.Range("L55").Value = 3750
$ws.Range("N55").Value = -4104
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 9000
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("M63").Value = -8251
$ws.Range("H64").Value = 962.5
$ws.Range("I64").Value = 833.3333
$ws.Range("K64").Value = 2499.9999
$ws.Range("M64").Value = -2229.9999
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 27000
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("M66").Value = -23256
$ws.Range("H67").Value = 962.5
$ws.Range("I67").Value = 833.3333
$ws.Range("K67").Value = 2499.9999
$ws.Range("M67").Value = -1563.9999
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H109").Value = 5000
$ws.Range("H118").Value = 3433.3333
$ws.Range("I118").Value = 3433.3333
$ws.Range("K118").Value = 10299.9999
$ws.Range("M118").Value = -9056.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 29999
$ws.Range("J100").Value = 29999
$ws.Range("L100").Value = 29999
$ws.Range("N100").Value = -32163

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3850
$ws.Range("I4").Value = 6900
$ws.Range("K4").Value = 6900
$ws.Range("M4").Value = -6787
$ws.Range("H7").Value = 5810.4
$ws.Range("I7").Value = 5351
$ws.Range("K7").Value = 5351
$ws.Range("M7").Value = -5239
$ws.Range("H28").Value = 3850
$ws.Range("I28").Value = 6900
$ws.Range("K28").Value = 6900
$ws.Range("M28").Value = -6668
$ws.Range("H37").Value = 3850
$ws.Range("I37").Value = 6900
$ws.Range("K37").Value = 6900
$ws.Range("M37").Value = -6793
$ws.Range("H55").Value = 417
$ws.Range("I55").Value = 450.5
$ws.Range("J55").Value = 350
$ws.Range("K55").Value = 450.5
$ws.Range("L55").Value = 350
$ws.Range("M55").Value = -277.5
$ws.Range("N55").Value = -696
$ws.Range("H68").Value = 60000
$ws.Range("I68").Value = 60000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 60000
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -59251
$ws.Range("H71").Value = 60000
$ws.Range("I71").Value = 60000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 300000
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -296256
$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 1600
$ws.Range("K100").Value = 1600
$ws.Range("M100").Value = -1059
$ws.Range("H126").Value = 5810.4
$ws.Range("I126").Value = 5351
$ws.Range("K126").Value = 16053
$ws.Range("M126").Value = -13583
$ws.Range("H132").Value = 4080.25
$ws.Range("I132").Value = 4290.3335
$ws.Range("J132").Value = 3450
$ws.Range("K132").Value = 12871.0005
$ws.Range("L132").Value = 10350
$ws.Range("M132").Value = -10341.0005
$ws.Range("N132").Value = -15410

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 4600
$ws.Range("J39").Value = 4600
$ws.Range("L39").Value = 4600
$ws.Range("N39").Value = -5426
$ws.Range("H119").Value = 18000
$ws.Range("J119").Value = 18000
$ws.Range("L119").Value = 18000
$ws.Range("N119").Value = -27676
